# Insert a new weekly data row at row 266 (pushing the existing rows down by one)
# and populate it with the new record for "Haba" at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(266).Insert()

$ws.Range("A266").Value = 9
$ws.Range("B266").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C266").Value = "Metropolitana"
$ws.Range("D266").Value = 45119
$ws.Range("E266").Value = 13
$ws.Range("F266").Value = 100112026
$ws.Range("G266").Value = "Haba"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 70
$ws.Range("K266").Value = 16000
$ws.Range("L266").Value = 18000
$ws.Range("M266").Value = 17000
$ws.Range("N266").Value = "`$/saco 25 kilos"
$ws.Range("O266").Value = "Provincia de Limarí"
$ws.Range("P266").Value = 680
$ws.Range("Q266").Value = 25
$ws.Range("R266").Value = "Hortaliza"
